# Actualización automática 2025-07-24 16:50:08
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("H7").Value = 1389.6
$wsGrupo.Range("I7").Value = 232.2
$wsGrupo.Range("M16").Value = 5497.62

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F7").Value = 2905.88
$wsMensual.Range("F16").Value = 5704.18
$wsMensual.Range("F56").Value = 68508.08

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D7").Value = 1816.2
$wsCumplimiento.Range("E7").Value = 583.8
$wsCumplimiento.Range("F7").Value = 0.75675

$wsCumplimiento.Range("D8").Value = 353.56
$wsCumplimiento.Range("E8").Value = 646.4400000000001
$wsCumplimiento.Range("F8").Value = 0.35356

$wsCumplimiento.Range("D16").Value = 55766.9
$wsCumplimiento.Range("E16").Value = -3940.440000000002
$wsCumplimiento.Range("F16").Value = 1.076031432592541

$wsCumplimiento.Range("D19").Value = 68508.08
$wsCumplimiento.Range("E19").Value = 45198.37064517916
$wsCumplimiento.Range("F19").Value = 0.6024995029858016

# Column D narrowed slightly (matches autofit-style width reduction in source diff)
$wsCumplimiento.Columns.Item(4).ColumnWidth = 12.166666666666666
